$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3, 0.02226487226428731),
    @(2, 4, 0.1942197725494594),
    @(2, 5, 0.2625407690567556),
    @(2, 6, 0.5831919000491652),
    @(2, 7, 0.2793809808603456),
    @(2, 8, 0.4517105476455114),
    @(2, 9, 0.318431087807987),
    @(2, 10, 0.4494539272028533),
    @(2, 13, 10.84818221959324),
    @(2, 15, 1.385613698168555),
    @(3, 3, 0.01946368688227551),
    @(3, 4, 0.1949775260558084),
    @(3, 5, 0.2447303387886919),
    @(3, 6, 0.6114998318860358),
    @(3, 7, 0.2872547831137737),
    @(3, 8, 0.4632136331428001),
    @(3, 9, 0.3248150021420884),
    @(3, 10, 0.4051300180841508),
    @(3, 13, 9.512207204291258),
    @(3, 15, 1.425936122475548),
    @(4, 3, 0.01773476375910121),
    @(4, 4, 0.1956772522827066),
    @(4, 5, 0.2340505805380388),
    @(4, 6, 0.6300841584133643),
    @(4, 7, 0.2927764708682687),
    @(4, 8, 0.4708084202134373),
    @(4, 9, 0.3292969436077975),
    @(4, 10, 0.3782208815716785),
    @(4, 13, 8.688592963863641),
    @(4, 15, 1.453251554731835),
    @(5, 3, 0.01702803181346013),
    @(5, 4, 0.1960205567107778),
    @(5, 5, 0.229760960970502),
    @(5, 6, 0.6379550542237986),
    @(5, 7, 0.2951960149413466),
    @(5, 8, 0.4740356636906071),
    @(5, 9, 0.3312625328013539),
    @(5, 10, 0.367328698961586),
    @(5, 13, 8.352105464001397),
    @(5, 15, 1.465015707467018),
    @(6, 3, 0.01691055007075448),
    @(6, 4, 0.1960810523530512),
    @(6, 5, 0.2290523933242596),
    @(6, 6, 0.6392798438445801),
    @(6, 7, 0.2956079121081814),
    @(6, 8, 0.4745794941889869),
    @(6, 9, 0.3315972580898787),
    @(6, 10, 0.3655244052621072),
    @(6, 13, 8.296179561475412),
    @(6, 15, 1.467007065022116),
    @(7, 3, 0.0177252412545883),
    @(7, 4, 0.1956816477398888),
    @(7, 5, 0.2339924786881014),
    @(7, 6, 0.6301891094183425),
    @(7, 7, 0.2928084201866454),
    @(7, 8, 0.4708514101293702),
    @(7, 9, 0.3293228916749271),
    @(7, 10, 0.3780736924777557),
    @(7, 13, 8.684058477147573),
    @(7, 15, 1.453407661025764),
    @(8, 3, 0.02130093284080203),
    @(8, 4, 0.1944318903473459),
    @(8, 5, 0.256345436200931),
    @(8, 6, 0.5926998035474895),
    @(8, 7, 0.2819512146087746),
    @(8, 8, 0.4555655663145188),
    @(8, 9, 0.3205142481139092),
    @(8, 10, 0.4341053968389872),
    @(8, 13, 10.38821143663955),
    @(8, 15, 1.398979976293305),
    @(9, 3, 0.02823863289725637),
    @(9, 4, 0.1938797188295922),
    @(9, 5, 0.3023042043232351),
    @(9, 6, 0.5289498476475174),
    @(9, 7, 0.2662646868179053),
    @(9, 8, 0.4298738527928663),
    @(9, 9, 0.3078007008672401),
    @(9, 10, 0.5465838604959572),
    @(9, 13, 13.70514322334702),
    @(9, 15, 1.312995346068064),
    @(10, 3, 0.03328724380502024),
    @(10, 4, 0.1946879194204598),
    @(10, 5, 0.3375092451963155),
    @(10, 6, 0.4883700950012191),
    @(10, 7, 0.258367602787402),
    @(10, 8, 0.4136971627431905),
    @(10, 9, 0.3013778854286784),
    @(10, 10, 0.6310748923327196),
    @(10, 13, 16.12947485575597),
    @(10, 15, 1.263102542856842),
    @(11, 3, 0.03557284802597849),
    @(11, 4, 0.195332428891632),
    @(11, 5, 0.3538719585472307),
    @(11, 6, 0.4713346292938496),
    @(11, 7, 0.2556090768305381),
    @(11, 8, 0.4069430900491824),
    @(11, 9, 0.2991197182559588),
    @(11, 10, 0.6699791705659379),
    @(11, 13, 17.23029923280217),
    @(11, 15, 1.243426911844409),
    @(12, 3, 0.03643670308539981),
    @(12, 4, 0.1956174609620689),
    @(12, 5, 0.360121100839379),
    @(12, 6, 0.4650942477034796),
    @(12, 7, 0.2546883016878922),
    @(12, 8, 0.404474080211628),
    @(12, 9, 0.298362562173395),
    @(12, 10, 0.6847842435452662),
    @(12, 13, 17.64692037473128),
    @(12, 15, 1.236422304067531),
    @(13, 3, 0.03625073102236342),
    @(13, 4, 0.1955542329107089),
    @(13, 5, 0.3587728331281852),
    @(13, 6, 0.4664287672922676),
    @(13, 7, 0.2548810372646386),
    @(13, 8, 0.4050018582232298),
    @(13, 9, 0.2985212314993291),
    @(13, 10, 0.6815923747773525),
    @(13, 13, 17.55720320398649),
    @(13, 15, 1.237910841088137),
    @(14, 3, 0.03564395141593479),
    @(14, 4, 0.1953550506857624),
    @(14, 5, 0.3543850008359044),
    @(14, 6, 0.4708169757259348),
    @(14, 7, 0.255530819971284),
    @(14, 8, 0.4067381777297783),
    @(14, 9, 0.2990554485649852),
    @(14, 10, 0.6711956999728557),
    @(14, 13, 17.26457933517418),
    @(14, 15, 1.242841629593158),
    @(15, 3, 0.03527206409580685),
    @(15, 4, 0.1952384168356076),
    @(15, 5, 0.3517043111214235),
    @(15, 6, 0.4735324761345368),
    @(15, 7, 0.2559450726202073),
    @(15, 8, 0.4078133125189112),
    @(15, 9, 0.2993955056240338),
    @(15, 10, 0.6648370937499521),
    @(15, 13, 17.08530954486667),
    @(15, 15, 1.245920330420432),
    @(16, 3, 0.033137646365347),
    @(16, 4, 0.1946514678783728),
    @(16, 5, 0.3364471474463215),
    @(16, 6, 0.4895125550474191),
    @(16, 7, 0.2585649913281003),
    @(16, 8, 0.4141508670338538),
    @(16, 9, 0.3015390313166435),
    @(16, 10, 0.6285422982145121),
    @(16, 13, 16.05749720377156),
    @(16, 15, 1.264450177200672),
    @(17, 3, 0.0318253755029474),
    @(17, 4, 0.1943631077525225),
    @(17, 5, 0.3271784675847016),
    @(17, 6, 0.4996847073363284),
    @(17, 7, 0.2603885643911781),
    @(17, 8, 0.4181948131340221),
    @(17, 9, 0.3030257720317735),
    @(17, 10, 0.6064004653000268),
    @(17, 13, 15.42648697414222),
    @(17, 15, 1.276599635785175),
    @(18, 3, 0.03106955722915927),
    @(18, 4, 0.1942232364310001),
    @(18, 5, 0.3218799145112428),
    @(18, 6, 0.505669187999878),
    @(18, 7, 0.2615156881235876),
    @(18, 8, 0.420577609508932),
    @(18, 9, 0.3039432133559927),
    @(18, 10, 0.5937088542182494),
    @(18, 13, 15.06335329408637),
    @(18, 15, 1.283871307155067),
    @(19, 3, 0.03081347451558258),
    @(19, 4, 0.1941803071985504),
    @(19, 5, 0.3200914191615567),
    @(19, 6, 0.5077182250057781),
    @(19, 7, 0.2619106349767222),
    @(19, 8, 0.4213940935729781),
    @(19, 9, 0.30426446794079),
    @(19, 10, 0.5894190531731738),
    @(19, 13, 14.9403677571085),
    @(19, 15, 1.286381728642198),
    @(20, 3, 0.03196517656748199),
    @(20, 4, 0.1943911065535673),
    @(20, 5, 0.3281617453232997),
    @(20, 6, 0.4985879828094895),
    @(20, 7, 0.2601863150012775),
    @(20, 8, 0.4177584351150756),
    @(20, 9, 0.3028610391219146),
    @(20, 10, 0.6087529313215043),
    @(20, 13, 15.4936787303526),
    @(20, 15, 1.275276866827085),
    @(21, 3, 0.03582222269362489),
    @(21, 4, 0.195412433857868),
    @(21, 5, 0.3556723525297372),
    @(21, 6, 0.4695222911303709),
    @(21, 7, 0.2553365711517301),
    @(21, 8, 0.406225760501556),
    @(21, 9, 0.2988958569440285),
    @(21, 10, 0.6742474325743331),
    @(21, 13, 17.35053604639251),
    @(21, 15, 1.241381135215192),
    @(22, 3, 0.03833335759749446),
    @(22, 4, 0.1963193201149807),
    @(22, 5, 0.373962183177909),
    @(22, 6, 0.4517563482205986),
    @(22, 7, 0.2528905584068326),
    @(22, 8, 0.3992057643418292),
    @(22, 9, 0.2968765923085499),
    @(22, 10, 0.7174793926360792),
    @(22, 13, 18.56274586896438),
    @(22, 15, 1.221834276785557),
    @(23, 3, 0.03699402357742088),
    @(23, 4, 0.1958129941760092),
    @(23, 5, 0.3641711682796256),
    @(23, 6, 0.4611238784599188),
    @(23, 7, 0.2541285249779293),
    @(23, 8, 0.4029045881427606),
    @(23, 9, 0.2979011047864972),
    @(23, 10, 0.6943646690007483),
    @(23, 13, 17.91587162084102),
    @(23, 15, 1.232024427053432),
    @(24, 3, 0.03190197678875961),
    @(24, 4, 0.1943783676760518),
    @(24, 5, 0.3277171118150619),
    @(24, 6, 0.4990833876244238),
    @(24, 7, 0.2602775070272827),
    @(24, 8, 0.4179555414651759),
    @(24, 9, 0.3029353197833302),
    @(24, 10, 0.6076892634052342),
    @(24, 13, 15.4633024548009),
    @(24, 15, 1.275873998284965),
    @(25, 3, 0.02637009676713831),
    @(25, 4, 0.1938202035343579),
    @(25, 5, 0.2896294282097642),
    @(25, 6, 0.5451192713015836),
    @(25, 7, 0.2698866757339715),
    @(25, 8, 0.4363565486522702),
    @(25, 9, 0.3107377176148489),
    @(25, 10, 0.5158506026481291),
    @(25, 13, 12.81024619287518),
    @(25, 15, 1.333970502984414)
)

foreach ($item in $data) {
    $ws.Cells.Item($item[0], $item[1]).Value = $item[2]
}
